$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for specific rows, matching the repulled/recalculated data.
$updates = @{
    5  = -3
    6  = 3
    8  = 8
    10 = 2
    12 = 4
    14 = -2
    15 = -1
    16 = -7
    17 = 2
    18 = -4
    20 = -2
    21 = -5
    22 = -1
    23 = 1
    25 = 8
    26 = -1
    27 = 2
    28 = 3
    29 = 5
    30 = -1
    31 = 3
    32 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
